# Apply updated crypto price/volume values (Wed Jun  5 03:29:48 UTC 2024 refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.831.22"
$ws.Range("E2").Value = "  +2.32%  "
$ws.Range("D3").Value = "3.807.32"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'702.79"
$ws.Range("E5").Value = "  +11.45%  "
$ws.Range("D6").Value = "'173.23"
$ws.Range("E6").Value = "  +4.29%  "
$ws.Range("D7").Value = "3.807.42"
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("D11").Value = "'7.68"
$ws.Range("E11").Value = "  +12.75%  "
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("E13").Value = "  +3.91%  "
$ws.Range("D14").Value = "'36.18"
$ws.Range("E14").Value = "  +3.88%  "
$ws.Range("D15").Value = "4.447.18"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "3.830.81"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "70.846.62"
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("D18").Value = "'17.72"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "'7.22"
$ws.Range("E19").Value = "  +2.78%  "
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").Value = "'11.36"
$ws.Range("E21").Value = "  +19.10%  "
$ws.Range("D22").Value = "'480.00"
$ws.Range("E22").Value = "  +3.44%  "
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("D24").Value = "'83.71"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("D25").Value = "'0.0000145"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").Value = "'12.37"
$ws.Range("E26").Value = "  +2.34%  "
$ws.Range("D27").Value = "'2.17"
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("D28").Value = "'10.41"
$ws.Range("E28").Value = "  +3.05%  "
$ws.Range("D29").Value = "3.958.12"
$ws.Range("E29").Value = "  +0.82%  "
$ws.Range("D31").Value = "'3.12"
$ws.Range("E31").Value = "  +16.35%  "
$ws.Range("E32").Value = "  +1.87%  "
$ws.Range("D33").Value = "'7.50"
$ws.Range("E33").Value = "  +5.68%  "
$ws.Range("D34").Value = "'29.59"
$ws.Range("E34").Value = "  +3.59%  "
$ws.Range("D35").Value = "'0.180"
$ws.Range("E35").Value = "  +6.03%  "
$ws.Range("D36").Value = "'9.22"
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "3.756.69"
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("D40").Value = "'3.51"
$ws.Range("E40").Value = "  +5.99%  "
$ws.Range("D41").Value = "'5.98"
$ws.Range("E41").Value = "  +3.06%  "
$ws.Range("D42").Value = "'0.000336"
$ws.Range("E42").Value = "  +24.76%  "
$ws.Range("D43").Value = "'2.19"
$ws.Range("E43").Value = "  +12.69%  "
$ws.Range("D44").Value = "'0.965"
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'45.63"
$ws.Range("E47").Value = "  +6.09%  "
$ws.Range("D48").Value = "'49.21"
$ws.Range("E48").Value = "  +5.35%  "
$ws.Range("D49").Value = "'159.60"
$ws.Range("E49").Value = "  +1.17%  "
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("D51").Value = "'0.300"
$ws.Range("E51").Value = "  +1.37%  "
